$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 7.8
$ws.Range("I2").Value = 1.4
$ws.Range("K2").Value = 17.5
$ws.Range("P2").Value = 2.16
$ws.Range("Q2").Value = 1.51
$ws.Range("BH2").Value = "2026-02-24 05:57:15"

# Row 3
$ws.Range("F3").Value = 1.88
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 4.9
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 3.95
$ws.Range("BH3").Value = "2026-02-24 05:57:15"

# Row 4
$ws.Range("BH4").Value = "2026-02-24 05:57:15"

# Row 5
$ws.Range("Q5").Value = 2.3
$ws.Range("BH5").Value = "2026-02-24 05:57:15"

# Row 6
$ws.Range("BH6").Value = "2026-02-24 05:57:15"

# Row 7
$ws.Range("BH7").Value = "2026-02-24 05:57:15"

# Row 8
$ws.Range("BH8").Value = "2026-02-24 05:57:15"

# Row 9
$ws.Range("P9").Value = 1.69
$ws.Range("BH9").Value = "2026-02-24 05:57:15"

# Row 10
$ws.Range("BH10").Value = "2026-02-24 05:57:15"
